# Added cases where activity and category does not exist
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table2 ("Input / Exists / Does Not Exist") edits ---
# Old row 5 (Exercise: Bike | Bike | Exercise) -> (Exercise: Bike | <blank> | Exercise, Bike)
$ws.Range("I5").ClearContents()
$ws.Range("J5").Value = "Exercise, Bike"

# Old row 6 (Exercise: Run | <blank> | Exercise, Run) -> (Does Not Exist category row becomes [Blank] example)
$ws.Range("H6").Value = "[Blank]"
$ws.Range("J6").ClearContents()

# Old row 7 ([Blank] | <blank> | <blank>) is removed entirely -> shrink the table by one row
$ws.Range("H7").ClearContents()

# Resize Table2 from H2:J7 to H2:J6 (drop the now-empty trailing row)
$lo2 = $ws.ListObjects.Item("Table2")
$lo2.Resize($ws.Range("H2:J6"))

# --- View state: move selection to K7 (topLeftCell scroll is a cosmetic view setting) ---
$ws.Range("K7").Select()
